# Update "want to go" (想去人数) counts in column F for both the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.
# The two sheets list largely the same events, but "全部类型" has one
# extra row inserted earlier in the sheet, so the target rows differ
# by 2 for the later entries.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# Updates for "展览" sheet (row -> new F value)
$expoUpdates = @{
    2  = 230
    5  = 13811
    11 = 90
    13 = 530
    16 = 13857
    18 = 617
    19 = 14926
    21 = 8224
    31 = 4
    35 = 1
    38 = 213
    41 = 5069
}

foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Updates for "全部类型" sheet (row -> new F value)
$allUpdates = @{
    2  = 230
    5  = 13811
    11 = 90
    13 = 530
    16 = 13857
    18 = 617
    19 = 14926
    21 = 8224
    31 = 4
    37 = 1
    40 = 213
    43 = 5069
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
